# "Made some changes to math and added Jalils combinations"
# Fill in Jalil's logged hours / task notes in column E/F for the week of
# 2019-12-03 .. 2019-12-09, add a "Bad"-styled "end" marker row, and clear
# the stray date that had been left in A28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hours worked each day (column E) and what was done (column F)
$ws.Range("E12").Value = "4 hours"
$ws.Range("E13").Value = "2 hours"
$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = "made the base for the banner"

$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "task 3basse work"

# Clear the leftover date in A28
$ws.Range("A28").Value = ""

# Mark the end of the tracked range with the built-in "Bad" cell style
$ws.Range("B27:I27").Style = "Bad"
$ws.Range("F27").Value = "end"

# Leave the selection where the last edit happened
$ws.Range("E19").Select()
